$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 0.9327180547071353

$ws.Range("D2:D6").Value = $newValue
